$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K1").Value = "Collection Type"
